$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp footer (row 1, column A)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 27 de Junio de 2020 a las 05:58"

# Refresh the country statistics: the data is a snapshot re-sorted in descending
# order by "Casos totales" (column B). Most rows keep the same rank, but some
# countries changed rank due to updated counts, so those rows' country name
# (column A) and/or numeric columns (B-H) need to be corrected to match the
# newly sorted snapshot.

$ws.Cells.Item(25, 2).Value = 83483
$ws.Cells.Item(25, 3).Value = 21
$ws.Cells.Item(25, 4).Value = 78444
$ws.Cells.Item(25, 5).Value = 405
$ws.Cells.Item(54, 2).Value = 20319
$ws.Cells.Item(54, 3).Value = 569
$ws.Cells.Item(54, 5).Value = 7621
$ws.Cells.Item(57, 1).Value = "Honduras"
$ws.Cells.Item(57, 2).Value = 15994
$ws.Cells.Item(57, 3).Value = 628
$ws.Cells.Item(57, 4).Value = 1678
$ws.Cells.Item(57, 5).Value = 13845
$ws.Cells.Item(57, 7).Value = 45
$ws.Cells.Item(57, 8).Value = 471
$ws.Cells.Item(58, 1).Value = "Ghana"
$ws.Cells.Item(58, 2).Value = 15834
$ws.Cells.Item(58, 3).Value = 0
$ws.Cells.Item(58, 4).Value = 11755
$ws.Cells.Item(58, 5).Value = 3976
$ws.Cells.Item(58, 7).Value = 0
$ws.Cells.Item(58, 8).Value = 103
$ws.Cells.Item(59, 1).Value = "Guatemala"
$ws.Cells.Item(59, 2).Value = 15828
$ws.Cells.Item(59, 3).Value = 209
$ws.Cells.Item(59, 4).Value = 3028
$ws.Cells.Item(59, 5).Value = 12128
$ws.Cells.Item(59, 7).Value = 49
$ws.Cells.Item(59, 8).Value = 672
$ws.Cells.Item(60, 1).Value = "Moldavia"
$ws.Cells.Item(60, 2).Value = 15776
$ws.Cells.Item(60, 4).Value = 8765
$ws.Cells.Item(60, 5).Value = 6496
$ws.Cells.Item(60, 8).Value = 515
$ws.Cells.Item(61, 1).Value = "Azerbaiyan"
$ws.Cells.Item(61, 2).Value = 15369
$ws.Cells.Item(61, 4).Value = 8364
$ws.Cells.Item(61, 5).Value = 6818
$ws.Cells.Item(61, 8).Value = 187
$ws.Cells.Item(81, 2).Value = 5722
$ws.Cells.Item(81, 3).Value = 179
$ws.Cells.Item(81, 4).Value = 641
$ws.Cells.Item(81, 5).Value = 4983
$ws.Cells.Item(81, 7).Value = 2
$ws.Cells.Item(81, 8).Value = 98
$ws.Cells.Item(87, 1).Value = "Venezuela"
$ws.Cells.Item(87, 2).Value = 4779
$ws.Cells.Item(87, 4).Value = 1327
$ws.Cells.Item(87, 5).Value = 3411
$ws.Cells.Item(87, 8).Value = 41
$ws.Cells.Item(88, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(88, 2).Value = 4643
$ws.Cells.Item(88, 4).Value = 4348
$ws.Cells.Item(88, 5).Value = 243
$ws.Cells.Item(88, 8).Value = 52
$ws.Cells.Item(165, 4).Value = 175
$ws.Cells.Item(165, 5).Value = 44
$ws.Cells.Item(176, 2).Value = 139
$ws.Cells.Item(176, 3).Value = 9
$ws.Cells.Item(176, 4).Value = 129
$ws.Cells.Item(176, 5).Value = 10
$ws.Cells.Item(201, 1).Value = "Laos"
$ws.Cells.Item(202, 1).Value = "Santa Lucia"
$ws.Cells.Item(203, 1).Value = "Dominica"
$ws.Cells.Item(204, 1).Value = "Fiyi"
$ws.Cells.Item(208, 1).Value = "Islas Malvinas"
$ws.Cells.Item(209, 1).Value = "Groenlandia"
$ws.Cells.Item(212, 1).Value = "Montserrat"
$ws.Cells.Item(212, 4).Value = 10
$ws.Cells.Item(212, 8).Value = 1
$ws.Cells.Item(213, 1).Value = "Seychelles"
$ws.Cells.Item(213, 4).Value = 11
$ws.Cells.Item(213, 8).Value = 0
